$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# URL changed to the new production / unilevel registration link
$ws.Range("B2").Value = "https://prod-unilevel.epixel.link/en/register/"

# Sponsor now lists the unilevel business-admin plus the two extra users
$ws.Range("B3").Value = "mpfp-base-unilevel-business-admin,user1,user2"

# Subdomain updated, with a note in C6 that it may need changing
$ws.Range("B6").Value = "antp087123"
$ws.Range("C6").Value = "If needed change"

# Enrollment Package simplified to "Bronze"
$ws.Range("B11").Value = "Bronze"

# The trailing blank marker row moves from row 19 down to row 23
$ws.Range("B19:AA19").Cut($ws.Range("B23"))
